$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MENDAFTAR")

# New peserta (participant) rows appended to the registration list.
$data = @(
    @("14211-3061", "Nita Febrina Butar-butar", "087744750232", 42187),
    @("14211-3771", "Sisilya Dewi Siregar",     "082276876382", 42188),
    @("13211-2964", "Indah Priskila Butar-butar","081281128784", 42188),
    @("13211-2476", "Imelda Sadna Sianturi",    "082370246730", 42188)
)

$startRow = 18
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $entry = $data[$i]

    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 4).Value = $entry[2]
    $ws.Cells.Item($row, 5).Value = $entry[3]
}

$ws.Range("D22").Select()
